$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell A2 text: append "[параметры лексемы]" before the trailing classifier bracket
$ws.Range("A2").Value = "синтактическая_категория [параметры лексемы] [классификатор_словоизменения]"

# Update the selection to A2 (was C2)
$ws.Range("A2").Select()
